# Update the interactive system comparison modules (TEA/LCA/MCDA) results
# Reflect a 1:0 outcome (sysA wins) instead of the previous 1:1 tie.

$wb = $excel.ActiveWorkbook

# --- Winner sheet ---
$wsWinner = $wb.Worksheets.Item("Winner")
$wsWinner.Range("B2").Value = "1:0"
$wsWinner.Range("D2").Value = "sysA"

# --- Score sheet ---
$wsScore = $wb.Worksheets.Item("Score")
$wsScore.Range("B2").Value = "1:0"
$wsScore.Range("D2").Value = 1
$wsScore.Range("E2").Value = 0

# --- Rank sheet ---
$wsRank = $wb.Worksheets.Item("Rank")
$wsRank.Range("B2").Value = "1:0"
$wsRank.Range("D2").Value = 1
$wsRank.Range("E2").Value = 2
